$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "is_locked_lbl" (D1) and "is_enabled_lbl" (E1) columns entirely.
# This shifts the former F1 (order_by) -> D1 and G1 (rem) -> E1,
# leaving the row with only 5 populated cells (A1:E1).
$ws.Range("D1:E1").EntireColumn.Delete() | Out-Null
